# Updated symbol list on Sun Jan 15 20:58:29 UTC 2023 with GitHub Actions
#
# Applies the cell-level text changes described by the diff against
# cryptos.xlsx / Sheet1. Columns B (Coin) and C (Link) are plain text and
# can be assigned directly. Columns D (Price) and E (Volume(1h)) hold
# numeric-/percent-looking strings (e.g. "301.50", "-0.69%") that must stay
# literal text (matching the original inlineStr cells) rather than being
# auto-converted by Excel into numbers/percentages. We force that by
# stamping the cell as Text ("@") before the write and then restoring the
# "Normal" cell style afterwards so no stray formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

function Set-PlainValue($cellRef, $text) {
    $ws.Range($cellRef).Value = $text
}

# Row 2 - BNB
Set-TextValue "D2" "301.50"
Set-TextValue "E2" "-0.69%"

# Row 3 - OKB
Set-TextValue "D3" "31.35"
Set-TextValue "E3" "-2.04%"

# Row 4 - HuobiToken
Set-TextValue "E4" "-3.05%"

# Row 5 - Cronos
Set-TextValue "D5" "0.07372"
Set-TextValue "E5" "-1.79%"

# Row 6 - FTXToken
Set-TextValue "D6" "2.430"
Set-TextValue "E6" "64.35%"

# Row 7 - KuCoinToken
Set-TextValue "E7" "1.11%"

# Row 8 (was MXToken, now GateToken)
Set-PlainValue "B8" "GateToken"
Set-PlainValue "C8" "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue "D8" "3.785"
Set-TextValue "E8" "-0.89%"

# Row 9 (was WazirX, now MXToken)
Set-PlainValue "B9" "MXToken"
Set-PlainValue "C9" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D9" "0.9163"
Set-TextValue "E9" "-0.39%"

# Row 10 (was LiechtensteinCryptoassetsExchange, now WazirX)
Set-PlainValue "B10" "WazirX"
Set-PlainValue "C10" "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue "D10" "0.1707"
Set-TextValue "E10" "0.50%"

# Row 11 (was MandalaExchangeToken, now LiechtensteinCryptoassetsExchange)
Set-PlainValue "B11" "LiechtensteinCryptoassetsExchange"
Set-PlainValue "C11" "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue "D11" "0.07607"
Set-TextValue "E11" "-3.19%"

# Row 12 (was BitrueCoin, now MandalaExchangeToken)
Set-PlainValue "B12" "MandalaExchangeToken"
Set-PlainValue "C12" "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue "D12" "0.08107"
Set-TextValue "E12" "1.04%"

# Row 13 (was BitMartToken, now BitrueCoin)
Set-PlainValue "B13" "BitrueCoin"
Set-PlainValue "C13" "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue "D13" "0.03038"
Set-TextValue "E13" "-0.44%"

# Row 14 (was BitForexToken, now BitMartToken)
Set-PlainValue "B14" "BitMartToken"
Set-PlainValue "C14" "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue "D14" "0.09923"
Set-TextValue "E14" "0.21%"

# Row 15 (was TigerCash, now BitForexToken)
Set-PlainValue "B15" "BitForexToken"
Set-PlainValue "C15" "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue "D15" "0.001505"
Set-TextValue "E15" "0.94%"

# Row 16 (was LEO, now TigerCash)
Set-PlainValue "B16" "TigerCash"
Set-PlainValue "C16" "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue "D16" "0.006167"
Set-TextValue "E16" "-0.12%"

# Row 17 (was GateToken, now LEO)
Set-PlainValue "B17" "LEO"
Set-PlainValue "C17" "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue "D17" "3.464"
Set-TextValue "E17" "0.14%"

# Row 19 - BitpandaEcosystemToken
Set-TextValue "D19" "0.3296"
Set-TextValue "E19" "-0.22%"

# Row 20 - ProBitToken
Set-TextValue "E20" "-0.51%"

# Row 21 - MCDex
Set-TextValue "D21" "4.646"
Set-TextValue "E21" "3.49%"

# Row 22 - CoinExToken
Set-TextValue "E22" "0.90%"

# Row 23 - ZBToken
Set-TextValue "D23" "0.1565"
Set-TextValue "E23" "-3.34%"

# Row 25 - HotbitToken
Set-TextValue "D25" "0.004486"
Set-TextValue "E25" "0.86%"

# Row 26 - NitroEx
Set-TextValue "D26" "0.0001297"

# Row 27 - UpBots
Set-TextValue "E27" "-3.16%"

# Row 39 - One
Set-TextValue "D39" "0.01736"
Set-TextValue "E39" "0.97%"

# Row 40 - IDEX
Set-TextValue "D40" "0.04523"
Set-TextValue "E40" "0.55%"

# Row 41 - KickToken
Set-TextValue "D41" "0.007208"
Set-TextValue "E41" "4.74%"

# Row 42 - BKEXToken
Set-TextValue "D42" "0.1346"
Set-TextValue "E42" "-0.23%"

# Row 43 - CEJI
Set-TextValue "E43" "1.34%"

# Row 44 - LocalTraders
Set-TextValue "D44" "0.01074"
Set-TextValue "E44" "-16.55%"

# Row 45 - CoinLion
Set-TextValue "D45" "0.00006271"
Set-TextValue "E45" "1.64%"

# Row 46 - CoinbaseStockToken
Set-TextValue "E46" "-33.26%"

# Row 47 - BOLO
Set-TextValue "D47" "0.8085"
Set-TextValue "E47" "13.69%"
